# Build site at 2022-09-26 16:07:08 UTC
# Rework of the "Docentes responsaveis" / "Programa" / "Avaliacao" block of
# the LOQ4008 syllabus sheet: several long-form paragraphs are swapped out
# for short values, a couple of rows collapse onto each other, and the
# trailing duplicated "Requisitos" row is dropped (24 rows -> 23 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ("Objetivos:") -------------------------------------------------
# Keep the label, replace the long objectives paragraph with the teacher
# name (re-used later on row 18 too).
$ws.Range("B10").Value = "849935 - Humberto Felipe da Silva"
$ws.Range("C10").Value = "849935 - Humberto Felipe da Silva"

# --- Row 13 ("Programa resumido:" / "Semestral") ---------------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14 ("Short syllabus:") --------------------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# --- Row 15 ("Programa:" / date) --------------------------------------------
$ws.Range("A15").Value = "Programa:"
# Copy the already-text "01/01/1996" cells from row 8 so the value lands as
# a shared string (not an auto-converted date serial) with the right style.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16 ("Syllabus:") ---------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# --- Row 17 ("Avaliação:") --------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# --- Row 18 ("Método:" / teacher name again) --------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B10").Copy($ws.Range("B18"))
$ws.Range("C10").Copy($ws.Range("C18"))
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19 ("Critério:" / written evaluation text) -------------------------
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Avaliação da disciplina constará de uma avaliação escrita programa e de um seminário  a ser apresentado pelos alunos no final do semestre."
$ws.Range("C19").Value = "Avaliação da disciplina constará de uma avaliação escrita programa e de um seminário  a ser apresentado pelos alunos no final do semestre."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20 ("Norma de recuperação:" / P1+P2 formula text) ------------------
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A avaliação escrita programa = P1 Seminário = P2   MP =(P1+P2)/2."
$ws.Range("C20").Value = "A avaliação escrita programa = P1 Seminário = P2   MP =(P1+P2)/2."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21 ("Bibliografia:" / recovery-exam paragraph) ---------------------
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Na recuperação haverá uma aula de revisão e na semana seguinte uma avaliação escrita. A média final será a média simples entre MP e nota da recuperação."
$ws.Range("C21").Value = "Na recuperação haverá uma aula de revisão e na semana seguinte uma avaliação escrita. A média final será a média simples entre MP e nota da recuperação."
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22 ("Requisitos:") --------------------------------------------------
# Only the label remains; the long bibliography paragraph that used to live
# here is dropped entirely.
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows.Item(22).AutoFit()

# --- Row 23: the prerequisite note moves up from row 24 ----------------------
$ws.Range("A23").Clear()
$ws.Range("B24").Copy($ws.Range("B23"))
$ws.Range("C24").Copy($ws.Range("C23"))
$ws.Rows.Item(23).RowHeight = 30

# Old row 24 now only duplicates row 23 - delete it so the sheet shrinks
# back down to 23 rows (A1:C23).
$ws.Rows.Item(24).Delete()
